# Add reading for next week: update row 11 of the schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: new "Topic" reading added for the following week (MAT212 schedule).
$ws.Range("A11").Value = "x"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = "T"
$ws.Range("D11").Value = "1/20/2025"
$ws.Range("E11").Value = "Logistic Regression: Inference"
$ws.Range("F11").Value = "/prepare/prep-10.qmd"
$ws.Range("G11").Value = " "
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

# Update the active cell selection to match the author's last edit location.
$ws.Range("I11").Select()
